$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three "ECs" sending-cluster rows (old rows 2-4); this shifts the
# "FAPs" sending-cluster rows (old 5-7) up to become rows 2-4.
$ws.Rows("2:4").Delete()

# Recalculated values (TPM update) for the remaining FAPs rows now in rows 2-4.
$row2 = @("FAPs","Wnt1","Fzd3","ECs",3,1,0.3532066666666667,1.05962,1,1,3,1,0.196431,0.589293,0.09717285149889213,0.09717285149889213,0.06938073874,0.62442664866,0.09717285149889213,0.09717285149889213)
$row3 = @("FAPs","Wnt1","Fzd3","FAPs",3,1,0.3532066666666667,1.05962,1,1,3,1,0.4307096666666667,1.292129,0.2130686423127578,0.2130686423127578,0.1521295256644445,1.36916573098,0.2130686423127578,0.2130686423127578)
$row4 = @("FAPs","Wnt1","Fzd3","MuSCs",3,1,0.3532066666666667,1.05962,1,1,3,1,1.394319,4.182957,0.68975850618835,0.68975850618835,0.49248276626,4.43234489634,0.68975850618835,0.68975850618835)

for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(2, $col).Value = $row2[$col - 1]
    $ws.Cells.Item(3, $col).Value = $row3[$col - 1]
    $ws.Cells.Item(4, $col).Value = $row4[$col - 1]
}
